$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the confidential note's date from 2021-04-09 to 2021-04-21
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-21 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2537561975349273
$ws.Range("E2").Value = 0.01094950155254115

$ws.Range("D3").Value = 0.4914754494355139
$ws.Range("E3").Value = 0.01133297355639473

$ws.Range("D4").Value = 0.1003716771038915
$ws.Range("E4").Value = 0.01319059453259408

$ws.Range("D5").Value = 0.09857351203003993
$ws.Range("E5").Value = 0.01607598525659171

$ws.Range("D6").Value = 0.05582316389562723
$ws.Range("E6").Value = 0.02108433734939763

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = 0.01243400499175085
